$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing all existing data down by one row
$ws.Rows.Item(1).Insert()

# Add header label for the data column
$ws.Range("A1").Value = "Data"

# Select the new header cell, mirroring the author's final selection state
$ws.Range("A2").Select()
